# Auto-generated edit script: updates the leads_output worksheet (rows 7-32).
# Commit: lead_scraper.py now prompts for the output filename and saves
# the collected leads inside a finally block, so partial results persist
# even on crash/interrupt. That changed the order/content of scraped rows
# and appended newly found leads to the bottom of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('A7').Value = 'Modern Letter R Real Estate Logo'
$ws.Range('B7').Value = '+1(215)796-4570'
$ws.Range('C7').Value = 'https://www.instagram.com/p/DPb03GeEgJK/'
$ws.Range('D7').Value = 'professionallogolads@gmail.com'
$ws.Range('A8').Value = 'Modern Luxury Interior Project At Bashundhara R/A + ...'
$ws.Range('B8').NumberFormat = '@'
$ws.Range('B8').Value = '0140339995'
$ws.Range('C8').Value = 'https://www.instagram.com/reel/DTFRT_FAeoA/'
$ws.Range('D8').Value = 'premiumtouchdecor@gmail.com'
$ws.Range('A10').Value = 'Peace Innovative Real Estate Ltd. (@peaceinnovative. ...'
$ws.Range('C10').Value = 'https://www.instagram.com/peaceinnovative.realestate/'
$ws.Range('D10').Value = 'peaceinnovative.realestate@gmail.com'
$ws.Range('A11').Value = 'Luxurious Apartment For Sale in Uttora (Sector 4, Road#6/ ...'
$ws.Range('C11').Value = 'https://www.instagram.com/p/CwAWoIMp3J9/'
$ws.Range('D11').Value = 'careers.studiodhakaltd@gmail.com'
$ws.Range('A12').Value = 'AIM Properties — Professional. Personal. Reliable. Finding ...'
$ws.Range('C12').Value = 'https://www.instagram.com/p/DR4FoafjWrz/'
$ws.Range('D12').Value = 'pmgroupebd@gmail.comread'
$ws.Range('A13').Value = 'ঢাকা ওয়েস্টার্ন ভ্যালি প্রকল্প এলাকায় নির্মাণ কাজ ...'
$ws.Range('C13').Value = 'https://www.instagram.com/reel/DTpyYe0DYBu/'
$ws.Range('D13').Value = 'info.faithrealestate@gmail.com'
$ws.Range('A14').Value = 'Luxury Living Redefined in Baridhara! 🏢✨🇧🇩 Concord ...'
$ws.Range('C14').Value = 'https://www.instagram.com/p/DTvlM5xlV-m/'
$ws.Range('D14').Value = 'pmgroupebd@gmail.com'
$ws.Range('A15').Value = 'Project : Troyee Number of Floors : G + 7 Living Floors Ground ...'
$ws.Range('C15').Value = 'https://www.instagram.com/reel/DUAqnqeglpd/'
$ws.Range('D15').Value = 'greenhutshafin@gmail.com'
$ws.Range('A16').Value = '#broucer'
$ws.Range('C16').Value = 'https://www.instagram.com/p/DU0joeKCB5x/'
$ws.Range('D16').Value = 'graphicsdac@gmail.com'
$ws.Range('A17').Value = '📢 Internship Opportunity – Night Shift (US Real Estate) ...'
$ws.Range('C17').Value = 'https://www.instagram.com/p/DTB849qk3PP/'
$ws.Range('D17').Value = 'bscgroupbd@gmail.com'
$ws.Range('A18').Value = ': Modern Flat for Sale in Dhaka! Looking for your dream home ...'
$ws.Range('C18').Value = 'https://www.instagram.com/reel/DUdpRqKjpy8/'
$ws.Range('D18').Value = 'dproperty.com.bd@gmail.com'
$ws.Range('A19').Value = 'New Launch Luxury Apartments! "REGEYA VILLA ...'
$ws.Range('C19').Value = 'https://www.instagram.com/p/DNd3kpSg4Yl/'
$ws.Range('D19').Value = 'structureshine@gmail.com'
$ws.Range('A20').Value = 'House No 59/A, Road No: 12/A, Dhanmondi, DHaka-1209. ...'
$ws.Range('C20').Value = 'https://www.instagram.com/p/CegXiQNv1A0/'
$ws.Range('D20').Value = 'studiodhakatld@gmail.com'
$ws.Range('A21').Value = 'Pixles Lab (@pixles_lab) · Dhaka, Bangladesh'
$ws.Range('C21').Value = 'https://www.instagram.com/pixles_lab/'
$ws.Range('D21').Value = 'thepixleslab@gmail.com'
$ws.Range('A22').Value = 'At REHAB Fair, Navana Real Estate presents a refined ...'
$ws.Range('C22').Value = 'https://www.instagram.com/reel/DSkJHTmkmqn/'
$ws.Range('D22').Value = 'arushrealtor@gmail.com'
$ws.Range('A23').Value = 'Ariyan Islam Rifat (@_mr.ariyan)'
$ws.Range('C23').Value = 'https://www.instagram.com/_mr.ariyan/'
$ws.Range('D23').Value = 'mr.ariyanislam2500@gmail.com'
$ws.Range('A24').Value = 'Welcome to JBS Gazi Landmark, a premium residential ...'
$ws.Range('C24').Value = 'https://www.instagram.com/reel/DTrnUR8D0Cq/'
$ws.Range('D24').Value = 'jsbuildersltdofficial@gmail.com'
$ws.Range('E24').Value = 'Google Search'
$ws.Range('A25').Value = '️ বাড়ি নির্মাণের আগে অবহেলা নয়— সঠিক সিদ্ধান্তই নিরাপদ ...'
$ws.Range('C25').Value = 'https://www.instagram.com/p/DUkYO8ukeii/'
$ws.Range('D25').Value = 'pisconsultantfimbd@gmail.com'
$ws.Range('E25').Value = 'Google Search'
$ws.Range('A26').Value = '🚧 Dhaka Western Valley Project – উন্নয়নের পথে আরও এক ধাপ ...'
$ws.Range('C26').Value = 'https://www.instagram.com/p/DTe3e4rAcAP/'
$ws.Range('D26').Value = 'dhakawesternvalley@gmail.com'
$ws.Range('E26').Value = 'Google Search'
$ws.Range('A27').Value = 'Real Estate Logo Design'
$ws.Range('C27').Value = 'https://www.instagram.com/p/C9k00kVyee3/?img_index=5'
$ws.Range('D27').Value = 'nusratnahianr@gmail.com'
$ws.Range('E27').Value = 'Google Search'
$ws.Range('A28').Value = 'IS THIS REALLY DHAKA? You have to see it to believe it!!'
$ws.Range('C28').Value = 'https://www.instagram.com/reel/DNqUDx502F-/'
$ws.Range('D28').Value = 'snowwspaces@gmail.com'
$ws.Range('E28').Value = 'Google Search'
$ws.Range('A29').Value = 'Anam & Naher Real Estate Ltd.(ANREL) (@anamnaher. ...'
$ws.Range('C29').Value = 'https://www.instagram.com/anamnaher.realestate/'
$ws.Range('E29').Value = 'Google Search'
$ws.Range('A30').Value = 'Notun Thikana (@notunthikana22) · Dhaka'
$ws.Range('C30').Value = 'https://www.instagram.com/notunthikana22/'
$ws.Range('E30').Value = 'Google Search'
$ws.Range('A31').Value = 'Purbachal Estate Agency (@estateagency.com.bd)'
$ws.Range('C31').Value = 'https://www.instagram.com/estateagency.com.bd/'
$ws.Range('E31').Value = 'Google Search'
$ws.Range('A32').Value = 'Maruf Raihan.Works (@marufraihan.works)'
$ws.Range('C32').Value = 'https://www-fallback.instagram.com/marufraihan.works/'
$ws.Range('E32').Value = 'Google Search'
